# Fruta / hortaliza, semanal
# Insert two new daily price records at rows 221-222 (pushing the existing
# records for rows 221-313 down to rows 223-315), matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 221 - shifts old rows 221:313 down to 223:315
$ws.Rows("221:222").Insert()

# --- New row 221 -----------------------------------------------------
$ws.Cells.Item(221, 1).Value2  = 4
$ws.Cells.Item(221, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(221, 3).Value2  = "Los Lagos"
$ws.Cells.Item(221, 4).Value2  = 44960
$ws.Cells.Item(221, 5).Value2  = 10
$ws.Cells.Item(221, 6).Value2  = "Fruta"
$ws.Cells.Item(221, 7).Value2  = 100108
$ws.Cells.Item(221, 8).Value2  = "Tropicales y subtropicales"
$ws.Cells.Item(221, 9).Value2  = 100108002
$ws.Cells.Item(221, 10).Value2 = "Mango"
$ws.Cells.Item(221, 11).Value2 = "Sin especificar"
$ws.Cells.Item(221, 12).Value2 = "Primera"
$ws.Cells.Item(221, 13).Value2 = 100
$ws.Cells.Item(221, 14).Value2 = 8000
$ws.Cells.Item(221, 15).Value2 = 8000
$ws.Cells.Item(221, 16).Value2 = 8000
$ws.Cells.Item(221, 17).Value2 = "`$/bandeja 4 kilos"
$ws.Cells.Item(221, 18).Value2 = "Perú"
$ws.Cells.Item(221, 19).Value2 = 2000
$ws.Cells.Item(221, 20).Value2 = 4

# --- New row 222 -----------------------------------------------------
$ws.Cells.Item(222, 1).Value2  = 4
$ws.Cells.Item(222, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(222, 3).Value2  = "Los Lagos"
$ws.Cells.Item(222, 4).Value2  = 44960
$ws.Cells.Item(222, 5).Value2  = 10
$ws.Cells.Item(222, 6).Value2  = "Fruta"
$ws.Cells.Item(222, 7).Value2  = 100108
$ws.Cells.Item(222, 8).Value2  = "Tropicales y subtropicales"
$ws.Cells.Item(222, 9).Value2  = 100108002
$ws.Cells.Item(222, 10).Value2 = "Mango"
$ws.Cells.Item(222, 11).Value2 = "Sin especificar"
$ws.Cells.Item(222, 12).Value2 = "Segunda"
$ws.Cells.Item(222, 13).Value2 = 100
$ws.Cells.Item(222, 14).Value2 = 7000
$ws.Cells.Item(222, 15).Value2 = 7000
$ws.Cells.Item(222, 16).Value2 = 7000
$ws.Cells.Item(222, 17).Value2 = "`$/bandeja 4 kilos"
$ws.Cells.Item(222, 18).Value2 = "Perú"
$ws.Cells.Item(222, 19).Value2 = 1750
$ws.Cells.Item(222, 20).Value2 = 4
